$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row: phone 79174409, blank birthday, 0 points.
$newRow = 10

# Phone numbers in column A are stored as text in this sheet's history of
# edits for "pending" rows (no birthday yet), so write it as text (leading
# apostrophe forces text interpretation) rather than a number.
$ws.Cells.Item($newRow, 1).Formula = "'79174409"
$ws.Cells.Item($newRow, 1).Style = $ws.Cells.Item(4, 1).Style

# Blank birthday cell (kept as an empty text cell, matching the other rows
# that have no birthday on file yet).
$ws.Cells.Item($newRow, 2).Formula = "'"
$ws.Cells.Item($newRow, 2).Style = $ws.Cells.Item(4, 2).Style

# Reset points to 0.00
$ws.Cells.Item($newRow, 3).Value = 0
